$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 344
$ws.Range("I28").Value = 252.58824
$ws.Range("K28").Value = 252.58824
$ws.Range("M28").Value = 232.41176
$ws.Range("H40").Value = 3180
$ws.Range("J40").Value = 4501
$ws.Range("L40").Value = 4501
$ws.Range("N40").Value = -4851
$ws.Range("H55").Value = 757.6667
$ws.Range("I55").Value = 519
$ws.Range("J55").Value = 1056
$ws.Range("K55").Value = 519
$ws.Range("L55").Value = 1056
$ws.Range("M55").Value = -305
$ws.Range("N55").Value = -1484
$ws.Range("H98").Value = 2984.3333
$ws.Range("I98").Value = 2981.625
$ws.Range("K98").Value = 2981.625
$ws.Range("M98").Value = -1483.625
$ws.Range("H107").Value = 1368.4736
$ws.Range("I107").Value = 1082.0588
$ws.Range("K107").Value = 1082.0588
$ws.Range("M107").Value = 837.9412
$ws.Range("H112").Value = 2843129.8
$ws.Range("I112").Value = 5499.5
$ws.Range("J112").Value = 3248505.5
$ws.Range("K112").Value = 16498.5
$ws.Range("L112").Value = 9745516.5
$ws.Range("M112").Value = -15390.5
$ws.Range("N112").Value = -9747732.5
$ws.Range("H118").Value = 242426
$ws.Range("I118").Value = 322833.16
$ws.Range("K118").Value = 968499.48
$ws.Range("M118").Value = -966842.48
$ws.Range("H122").Value = 2984.3333
$ws.Range("I122").Value = 2981.625
$ws.Range("K122").Value = 8944.875
$ws.Range("M122").Value = -6494.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1700.0625
$ws.Range("I45").Value = 1630.5
$ws.Range("K45").Value = 1630.5
$ws.Range("M45").Value = -1253.5
$ws.Range("H61").Value = 25654900
$ws.Range("I61").Value = 66678860
$ws.Range("K61").Value = 66678860
$ws.Range("M61").Value = -66678648
$ws.Range("H110").Value = 6098956
$ws.Range("I110").Value = 7576521
$ws.Range("K110").Value = 7576521
$ws.Range("M110").Value = -7574476
$ws.Range("H122").Value = 1106
$ws.Range("I122").Value = 1112
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 3336
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -886
$ws.Range("N122").Value = -8200
$ws.Range("H132").Value = 6067885
$ws.Range("I132").Value = 10005281
$ws.Range("K132").Value = 30015843
$ws.Range("M132").Value = -30013313
$ws.Range("H136").Value = 25654900
$ws.Range("I136").Value = 66678860
$ws.Range("K136").Value = 200036580
$ws.Range("M136").Value = -200034030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 75263.85000000001
$ws.Range("I86").Value = 1164.96
$ws.Range("K86").Value = 1164.96
$ws.Range("M86").Value = -41.96000000000004
$ws.Range("H89").Value = 75263.85000000001
$ws.Range("I89").Value = 1164.96
$ws.Range("K89").Value = 5824.8
$ws.Range("M89").Value = -208.8000000000002
$ws.Range("H107").Value = 1198.25
$ws.Range("I107").Value = 916.0741
$ws.Range("J107").Value = 2044.7778
$ws.Range("K107").Value = 916.0741
$ws.Range("L107").Value = 2044.7778
$ws.Range("M107").Value = 1003.9259
$ws.Range("N107").Value = -5884.7778
$ws.Range("H134").Value = 5227.625
$ws.Range("I134").Value = 3971.8
$ws.Range("J134").Value = 7320.6665
$ws.Range("K134").Value = 11915.4
$ws.Range("L134").Value = 21961.9995
$ws.Range("M134").Value = -9380.400000000001
$ws.Range("N134").Value = -27031.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 514.5
$ws.Range("I5").Value = 450
$ws.Range("J5").Value = 542.1429000000001
$ws.Range("K5").Value = 450
$ws.Range("L5").Value = 542.1429000000001
$ws.Range("M5").Value = -338
$ws.Range("N5").Value = -766.1429000000001
$ws.Range("H16").Value = 1708.0476
$ws.Range("I16").Value = 1435.0667
$ws.Range("K16").Value = 1435.0667
$ws.Range("M16").Value = -1148.0667
$ws.Range("H62").Value = 7308
$ws.Range("I62").Value = 6419
$ws.Range("J62").Value = 7974.75
$ws.Range("K62").Value = 6419
$ws.Range("L62").Value = 7974.75
$ws.Range("M62").Value = -5795
$ws.Range("N62").Value = -9222.75
$ws.Range("H65").Value = 7308
$ws.Range("I65").Value = 6419
$ws.Range("J65").Value = 7974.75
$ws.Range("K65").Value = 32095
$ws.Range("L65").Value = 39873.75
$ws.Range("M65").Value = -28975
$ws.Range("N65").Value = -46113.75
$ws.Range("H113").Value = 1708.0476
$ws.Range("I113").Value = 1435.0667
$ws.Range("K113").Value = 1435.0667
$ws.Range("M113").Value = 734.9332999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 564.55
$ws.Range("I34").Value = 67.789474
$ws.Range("K34").Value = 203.368422
$ws.Range("M34").Value = -119.368422
$ws.Range("H39").Value = 7906.25
$ws.Range("I39").Value = 250
$ws.Range("K39").Value = 750
$ws.Range("M39").Value = -456
$ws.Range("H55").Value = 2319.875
$ws.Range("J55").Value = 2970
$ws.Range("L55").Value = 8910
$ws.Range("N55").Value = -9264
$ws.Range("H56").Value = 11417
$ws.Range("I56").Value = 11417
$ws.Range("K56").Value = 11417
$ws.Range("M56").Value = -10887
$ws.Range("H62").Value = 2073.0146
$ws.Range("I62").Value = 1249.0278
$ws.Range("K62").Value = 3747.0834
$ws.Range("M62").Value = -3061.0834
$ws.Range("H65").Value = 2073.0146
$ws.Range("I65").Value = 1249.0278
$ws.Range("K65").Value = 11241.2502
$ws.Range("M65").Value = -7809.2502
$ws.Range("H82").Value = 10998.5
$ws.Range("I82").Value = 10998.5
$ws.Range("K82").Value = 32995.5
$ws.Range("M82").Value = -32589.5
$ws.Range("H85").Value = 10998.5
$ws.Range("I85").Value = 10998.5
$ws.Range("K85").Value = 32995.5
$ws.Range("M85").Value = -31591.5
$ws.Range("H130").Value = 4500
$ws.Range("J130").Value = 4500
$ws.Range("L130").Value = 13500
$ws.Range("N130").Value = -23540
$ws.Range("H137").Value = 3425.5557
$ws.Range("I137").Value = 2215
$ws.Range("J137").Value = 3771.4285
$ws.Range("K137").Value = 6645
$ws.Range("L137").Value = 11314.2855
$ws.Range("M137").Value = -1545
$ws.Range("N137").Value = -21514.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1874.6666
$ws.Range("I122").Value = 1999.6
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 5998.799999999999
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -3548.799999999999
$ws.Range("N122").Value = -8650
$ws.Range("H126").Value = 3366.875
$ws.Range("I126").Value = 2338.8572
$ws.Range("J126").Value = 4166.4443
$ws.Range("K126").Value = 7016.571599999999
$ws.Range("L126").Value = 12499.3329
$ws.Range("M126").Value = -4546.571599999999
$ws.Range("N126").Value = -17439.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3574.1538
$ws.Range("I22").Value = 2538.8
$ws.Range("J22").Value = 4221.25
$ws.Range("K22").Value = 2538.8
$ws.Range("L22").Value = 4221.25
$ws.Range("M22").Value = -2243.8
$ws.Range("N22").Value = -4811.25
$ws.Range("H27").Value = 3574.1538
$ws.Range("I27").Value = 2538.8
$ws.Range("J27").Value = 4221.25
$ws.Range("K27").Value = 2538.8
$ws.Range("L27").Value = 4221.25
$ws.Range("M27").Value = -2431.8
$ws.Range("N27").Value = -4435.25
$ws.Range("H46").Value = 6852.8623
$ws.Range("I46").Value = 1900.5
$ws.Range("J46").Value = 7219.7036
$ws.Range("K46").Value = 1900.5
$ws.Range("L46").Value = 7219.7036
$ws.Range("M46").Value = -1712.5
$ws.Range("N46").Value = -7595.7036
$ws.Range("H59").Value = 43069
$ws.Range("J59").Value = 43069
$ws.Range("L59").Value = 43069
$ws.Range("N59").Value = -44377
$ws.Range("H68").Value = 4749.75
$ws.Range("I68").Value = 4166.3335
$ws.Range("J68").Value = 6500
$ws.Range("K68").Value = 4166.3335
$ws.Range("L68").Value = 6500
$ws.Range("M68").Value = -3417.3335
$ws.Range("N68").Value = -7998
$ws.Range("H71").Value = 4749.75
$ws.Range("I71").Value = 4166.3335
$ws.Range("J71").Value = 6500
$ws.Range("K71").Value = 20831.6675
$ws.Range("L71").Value = 32500
$ws.Range("M71").Value = -17087.6675
$ws.Range("N71").Value = -39988
$ws.Range("H111").Value = 37369
$ws.Range("I111").Value = 24351
$ws.Range("J111").Value = 50387
$ws.Range("K111").Value = 24351
$ws.Range("L111").Value = 50387
$ws.Range("M111").Value = -20261
$ws.Range("N111").Value = -58567
$ws.Range("H123").Value = 78000
$ws.Range("J123").Value = 78000
$ws.Range("L123").Value = 78000
$ws.Range("N123").Value = -87800
$ws.Range("H124").Value = 67774.836
$ws.Range("J124").Value = 67774.836
$ws.Range("L124").Value = 67774.836
$ws.Range("N124").Value = -77594.836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1875.0714
$ws.Range("I6").Value = 2400
$ws.Range("J6").Value = 1481.375
$ws.Range("K6").Value = 2400
$ws.Range("L6").Value = 1481.375
$ws.Range("M6").Value = -2285
$ws.Range("N6").Value = -1711.375
$ws.Range("H95").Value = 59739.57
$ws.Range("J95").Value = 59739.57
$ws.Range("L95").Value = 59739.57
$ws.Range("N95").Value = -65231.57
$ws.Range("H107").Value = 1089.2778
$ws.Range("I107").Value = 724.2308
$ws.Range("J107").Value = 2038.4
$ws.Range("K107").Value = 2172.6924
$ws.Range("L107").Value = 6115.200000000001
$ws.Range("M107").Value = -252.6923999999999
$ws.Range("N107").Value = -9955.200000000001
$ws.Range("H127").Value = 24214
$ws.Range("J127").Value = 24214
$ws.Range("L127").Value = 24214
$ws.Range("N127").Value = -34134
$ws.Range("H132").Value = 8887.083000000001
$ws.Range("J132").Value = 9850
$ws.Range("L132").Value = 29550
$ws.Range("N132").Value = -34610
